$d = $word.ActiveDocument

# --- Portuguese "Programa" paragraph: split the single run into 8 runs
#     separated by manual line breaks (<w:br/>), one before each numbered item.

$d.Content.Find.Execute("oxigênio dissolvido.2) Determinação", $false, $false, $false, $false, $false, $true, 1, $false, "oxigênio dissolvido.^l2) Determinação", 2)
$d.Content.Find.Execute("e aeração.3) Determinação", $false, $false, $false, $false, $false, $true, 1, $false, "e aeração.^l3) Determinação", 2)
$d.Content.Find.Execute("de descoloração.4) Imobilização", $false, $false, $false, $false, $false, $true, 1, $false, "de descoloração.^l4) Imobilização", 2)
$d.Content.Find.Execute("células imobilizadas.5) Realização", $false, $false, $false, $false, $false, $true, 1, $false, "células imobilizadas.^l5) Realização", 2)
$d.Content.Find.Execute("substrato e produto.6) Projeto", $false, $false, $false, $false, $false, $true, 1, $false, "substrato e produto.^l6) Projeto", 2)
$d.Content.Find.Execute("ideais para aplicação.7) Caracterização", $false, $false, $false, $false, $false, $true, 1, $false, "ideais para aplicação.^l7) Caracterização", 2)
$d.Content.Find.Execute("atividade específica.8) Planejamento", $false, $false, $false, $false, $false, $true, 1, $false, "atividade específica.^l8) Planejamento", 2)

# --- English (italic) paragraph: same split.

$d.Content.Find.Execute("dissolved oxygen probes.2)Experimental", $false, $false, $false, $false, $false, $true, 1, $false, "dissolved oxygen probes.^l2)Experimental", 2)
$d.Content.Find.Execute("aeration conditions.3)Experimental", $false, $false, $false, $false, $false, $true, 1, $false, "aeration conditions.^l3)Experimental", 2)
$d.Content.Find.Execute("decolorization method.4)Immobilization", $false, $false, $false, $false, $false, $true, 1, $false, "decolorization method.^l4)Immobilization", 2)
$d.Content.Find.Execute("immobilized cells.5)Execution", $false, $false, $false, $false, $false, $true, 1, $false, "immobilized cells.^l5)Execution", 2)
$d.Content.Find.Execute("product concentrations.6)Enzyme", $false, $false, $false, $false, $false, $true, 1, $false, "product concentrations.^l6)Enzyme", 2)
$d.Content.Find.Execute("for application.7)Characterization", $false, $false, $false, $false, $false, $true, 1, $false, "for application.^l7)Characterization", 2)
$d.Content.Find.Execute("specific activity.8)Experimental", $false, $false, $false, $false, $false, $true, 1, $false, "specific activity.^l8)Experimental", 2)
